$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B8").Value = "TbLocalizationConfig"
$ws.Range("C8").Value = "LocalizationConfig"
$ws.Range("D7").Copy($ws.Range("D8"))
$ws.Range("E8").Value = "../LocalizationConfig@本地化.xlsx"
$ws.Range("G4").Copy($ws.Range("G8"))

$ws.Range("B9").Value = "TbEquipmentConfig"
$ws.Range("C9").Value = "EquipmentConfig"
$ws.Range("D7").Copy($ws.Range("D9"))
$ws.Range("E9").Value = "../EquipmentConfig@装备.xlsx"
$ws.Range("G4").Copy($ws.Range("G9"))
